$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows before row 153, pushing the existing rows 153-160 down
# to 156-163 (Excel copies formatting, e.g. the date style on column D, from
# the row above automatically).
$ws.Rows("153:155").Insert()

# Row 153: new weekly record (Angeleno / Especial)
$ws.Range("A153").Value2 = 7
$ws.Range("B153").Value2 = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C153").Value2 = "Ñuble"
$ws.Range("D153").Value2 = 45075
$ws.Range("E153").Value2 = 16
$ws.Range("F153").Value2 = "Fruta"
$ws.Range("G153").Value2 = 100103
$ws.Range("H153").Value2 = "Frutos de hueso (carozo)"
$ws.Range("I153").Value2 = 100103002
$ws.Range("J153").Value2 = "Ciruela"
$ws.Range("K153").Value2 = "Angeleno"
$ws.Range("L153").Value2 = "Especial"
$ws.Range("M153").Value2 = 40
$ws.Range("N153").Value2 = 10000
$ws.Range("O153").Value2 = 10000
$ws.Range("P153").Value2 = 10000
$ws.Range("Q153").Value2 = "`$/bandeja 18 kilos granel"
$ws.Range("R153").Value2 = "Región de O'Higgins"
$ws.Range("S153").Value2 = 556
$ws.Range("T153").Value2 = 18

# Row 154: new weekly record (Angeleno / Primera)
$ws.Range("A154").Value2 = 7
$ws.Range("B154").Value2 = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C154").Value2 = "Ñuble"
$ws.Range("D154").Value2 = 45075
$ws.Range("E154").Value2 = 16
$ws.Range("F154").Value2 = "Fruta"
$ws.Range("G154").Value2 = 100103
$ws.Range("H154").Value2 = "Frutos de hueso (carozo)"
$ws.Range("I154").Value2 = 100103002
$ws.Range("J154").Value2 = "Ciruela"
$ws.Range("K154").Value2 = "Angeleno"
$ws.Range("L154").Value2 = "Primera"
$ws.Range("M154").Value2 = 50
$ws.Range("N154").Value2 = 8000
$ws.Range("O154").Value2 = 8000
$ws.Range("P154").Value2 = 8000
$ws.Range("Q154").Value2 = "`$/bandeja 18 kilos granel"
$ws.Range("R154").Value2 = "Región de O'Higgins"
$ws.Range("S154").Value2 = 444
$ws.Range("T154").Value2 = 18

# Row 155: new weekly record (Angeleno / Segunda)
$ws.Range("A155").Value2 = 7
$ws.Range("B155").Value2 = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C155").Value2 = "Ñuble"
$ws.Range("D155").Value2 = 45075
$ws.Range("E155").Value2 = 16
$ws.Range("F155").Value2 = "Fruta"
$ws.Range("G155").Value2 = 100103
$ws.Range("H155").Value2 = "Frutos de hueso (carozo)"
$ws.Range("I155").Value2 = 100103002
$ws.Range("J155").Value2 = "Ciruela"
$ws.Range("K155").Value2 = "Angeleno"
$ws.Range("L155").Value2 = "Segunda"
$ws.Range("M155").Value2 = 30
$ws.Range("N155").Value2 = 6000
$ws.Range("O155").Value2 = 6000
$ws.Range("P155").Value2 = 6000
$ws.Range("Q155").Value2 = "`$/bandeja 18 kilos granel"
$ws.Range("R155").Value2 = "Región de O'Higgins"
$ws.Range("S155").Value2 = 333
$ws.Range("T155").Value2 = 18
